$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat is forced to text ("@") before assignment so that values which
# look numeric (e.g. "1.0000", "0.07753") are not silently coerced into numbers
# by Excel, and the Style is reset back to Normal afterwards so no formatting
# residue is left behind on the cell.
$cellValues = @{
    "D2" = "29.313.81"
    "E2" = "  +0.43%  "
    "D3" = "1.875.92"
    "D4" = "1.0000"
    "E4" = "  -0.01%  "
    "D5" = "0.7135"
    "E5" = "  +0.07%  "
    "D6" = "242.17"
    "E6" = "  +0.79%  "
    "E7" = "  -0.01%  "
    "D8" = "0.3108"
    "E8" = "  +1.16%  "
    "D9" = "0.07753"
    "E9" = "  +0.03%  "
    "D10" = "24.93"
    "E10" = "  -0.57%  "
    "D11" = "0.08539"
    "E11" = "  +3.50%  "
    "D12" = "1.883.83"
    "E12" = "  +1.67%  "
    "D13" = "5.219"
    "E13" = "  -0.20%  "
    "D14" = "0.7106"
    "E14" = "  -0.68%  "
    "D15" = "91.47"
    "E15" = "  +1.38%  "
    "D16" = "29.306.33"
    "E16" = "  +0.49%  "
    "D17" = "0.000008200"
    "E17" = "  +5.28%  "
    "D18" = "6.007"
    "E18" = "  +2.56%  "
    "D19" = "241.97"
    "E19" = "  -0.62%  "
    "D20" = "2.134.73"
    "E20" = "  +2.09%  "
    "D21" = "13.25"
    "E21" = "  +0.77%  "
    "D22" = "0.9996"
    "E22" = "  -0.03%  "
    "D23" = "7.805"
    "E24" = "  +0.01%  "
    "E25" = "  +0.68%  "
    "D26" = "162.88"
    "E26" = "  +0.30%  "
    "D27" = "9.050"
    "E27" = "  +1.64%  "
    "D28" = "18.49"
    "E28" = "  +1.01%  "
    "E29" = "  +1.23%  "
    "D30" = "4.399"
    "E30" = "  -0.04%  "
    "D31" = "4.319"
    "E31" = "  +2.81%  "
    "E32" = "  -2.51%  "
    "D33" = "0.05257"
    "E33" = "  +1.43%  "
    "D34" = "1.933"
    "E34" = "  +1.26%  "
    "E35" = "  +0.32%  "
    "D36" = "0.7454"
    "E36" = "  +2.70%  "
    "D37" = "2.685"
    "E37" = "  +0.32%  "
    "D38" = "0.01867"
    "E38" = "  +0.68%  "
    "E39" = "  +1.12%  "
    "D40" = "1.183.35"
    "E40" = "  +1.76%  "
    "D41" = "6.386"
    "E41" = "  +3.79%  "
    "D42" = "72.95"
    "E42" = "  +1.10%  "
    "D43" = "0.8868"
    "E43" = "  -1.80%  "
    "D44" = "106.39"
    "E44" = "  +4.89%  "
    "D45" = "1.0000"
    "D46" = "2.030.16"
    "E46" = "  +1.44%  "
    "D47" = "1.810"
    "E47" = "  +2.62%  "
    "D48" = "0.5207"
    "E48" = "  -0.16%  "
    "E49" = "  +1.52%  "
    "D50" = "9.392"
    "E50" = "  +0.89%  "
    "D51" = "0.4316"
    "E51" = "  +1.33%  "
}

foreach ($addr in $cellValues.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $cellValues[$addr]
    $range.Style = "Normal"
}
